# ============================================================
# cs-en-us-066pct.xlsx weekly refresh: shift the reporting week
# forward by one week (7/10-7/16/2023 -> 7/17-7/23/2023, Volume
# 30 Number 28 -> 29) and replace the crime-count table (rows
# 15-30) with the newly collected weekly figures.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings: volume/number + reporting week dates ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Style-donor cells (stable, untouched by this weeks refresh) ---
# I14 / K14 already carry the s=15 (#,##0) / s=16 (#,##0.0) numeric
# styles; D28 already carries the s=14 (General, text-placeholder)
# style used for "0"/"***.+" blanks. We PasteSpecial-copy just the
# formatting from these donors onto cells that flip between the
# text-placeholder style and a numeric style this week.

# --- Cells flipping from text-placeholder to numeric (new activity) ---
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("I14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
$ws.Range("I14").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("G26").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Value = -100
$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

# --- Cells flipping from numeric back to the text placeholder (no activity) ---
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D28").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# --- Remaining numeric updates (style/type unchanged) ---
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -16.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 46
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 2.222222222222
$ws.Range("L16").Value = 21.052631578947
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -88.697788697788
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 119
$ws.Range("J17").Value = 119
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 52.564102564102
$ws.Range("M17").Value = 41.666666666666
$ws.Range("N17").Value = -36.021505376344
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 65
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = -35
$ws.Range("L18").Value = -1.515151515151
$ws.Range("M18").Value = -69.194312796208
$ws.Range("N18").Value = -93.062966915688
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.769230769230
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -31.481481481481
$ws.Range("I19").Value = 301
$ws.Range("J19").Value = 330
$ws.Range("K19").Value = -8.787878787878
$ws.Range("L19").Value = 25.941422594142
$ws.Range("M19").Value = 38.073394495412
$ws.Range("N19").Value = -15.686274509803
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("I20").Value = 86
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 72
$ws.Range("L20").Value = 95.454545454545
$ws.Range("M20").Value = 7.5
$ws.Range("N20").Value = -91.863765373699
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -37.037037037037
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -18.811881188118
$ws.Range("I21").Value = 629
$ws.Range("J21").Value = 657
$ws.Range("K21").Value = -4.261796042617
$ws.Range("L21").Value = 33.829787234042
$ws.Range("M21").Value = -8.442503639010
$ws.Range("N21").Value = -78.771515356058
$ws.Range("M22").Value = -54.545454545454
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -5.263157894736
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = -7.5
$ws.Range("I24").Value = 627
$ws.Range("J24").Value = 620
$ws.Range("K24").Value = 1.129032258064
$ws.Range("L24").Value = 46.838407494145
$ws.Range("M24").Value = 27.439024390243
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 33
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 195
$ws.Range("J25").Value = 177
$ws.Range("K25").Value = 10.169491525423
$ws.Range("L25").Value = 23.417721518987
$ws.Range("M25").Value = -22.310756972111
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = -26.666666666666
$ws.Range("L26").Value = 37.5
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -45.454545454545
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = -21.276595744680
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = -44.444444444444
$ws.Range("L30").Value = -28.571428571428

$excel.CutCopyMode = 0
